$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52, pushing existing rows 52-108 down to 53-109.
$ws.Rows(52).Insert()

# Populate the newly inserted row 52 with the new weekly data point.
$ws.Range("A52").Value = 1
$ws.Range("B52").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C52").Value = "Arica y Parinacota"
$ws.Range("D52").NumberFormat = $ws.Range("D53").NumberFormat
$ws.Range("D52").Value = 45280
$ws.Range("E52").Value = 15
$ws.Range("F52").Value = 100112031
$ws.Range("G52").Value = "Poroto verde"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 1500
$ws.Range("K52").Value = 500
$ws.Range("L52").Value = 600
$ws.Range("M52").Value = 550
$ws.Range("N52").Value = "$/kilo"
$ws.Range("O52").Value = "Región de Arica y Parinacota"
$ws.Range("P52").Value = 550
$ws.Range("Q52").Value = 1
$ws.Range("R52").Value = "Hortaliza"
